$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("costs")

# Row 4: Project Y (cost)
$ws.Range("B4").Value = "Project Y"
$ws.Range("C4").Value = -100000

# New header cells for the PnL timing columns
$ws.Range("D2").Value = "PnL Start"
$ws.Range("E2").Value = "PnL End"

# Row 5: Project Z (income)
$ws.Range("B5").Value = "Project Z"
$ws.Range("C5").Value = 10000000

# Copy the existing date cell's format (style with the date number format)
# onto the new date cells so they match the workbook's existing date styling.
$ws.Range("A3").Copy() | Out-Null
$ws.Range("A4:A5").PasteSpecial(-4122) | Out-Null
$ws.Range("D4:E5").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("A4").Value = 46752
$ws.Range("A5").Value = 45991
$ws.Range("D4").Value = 46767
$ws.Range("E4").Value = 47026
$ws.Range("D5").Value = 45703
$ws.Range("E5").Value = 45793

# Match the new column widths used for the PnL timing columns
$ws.Columns.Item(4).ColumnWidth = 9.3
$ws.Columns.Item(5).ColumnWidth = 9.3

# Move the selection like in the authored workbook
[void]$ws.Range("C10").Select()
